$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of an existing data cell across the new
# K:S columns (rows 1-9) so the new cells pick up style index "1" just
# like the existing A:J columns, instead of defaulting to no style.
$ws.Range("A1").Copy()
$ws.Range("K1:S9").PasteSpecial(-4122)

# --- New Tiefling subrace columns (K:S) ------------------------------
$ws.Range("K1").Value = "Asmodeus Tiefling"
$ws.Range("L1").Value = "Baalzebul Tiefling"
$ws.Range("M1").Value = "Dispater Tiefling"
$ws.Range("N1").Value = "Fierna Tiefling"
$ws.Range("O1").Value = "Glasya Tiefling"
$ws.Range("P1").Value = "Levistus Tiefling"
$ws.Range("Q1").Value = "Mammon Tiefling"
$ws.Range("R1").Value = "Mephistopheles Tiefling"
$ws.Range("S1").Value = "Zariel Tiefling"
$ws.Range("K2").Value = "0=0=0=1=0=2"
$ws.Range("L2").Value = "0=0=0=1=0=2"
$ws.Range("M2").Value = "0=1=0=0=0=2"
$ws.Range("N2").Value = "0=0=0=0=1=2"
$ws.Range("O2").Value = "0=1=0=0=0=2"
$ws.Range("P2").Value = "0=0=1=0=0=2"
$ws.Range("Q2").Value = "0=0=0=1=0=2"
$ws.Range("R2").Value = "0=0=0=1=0=2"
$ws.Range("S2").Value = "1=0=0=0=0=2"
$ws.Range("K3").Value = "Charisma=1/Thaumaturgy/0=3/Hellish Rebuke/2=5/Darkness/2"
$ws.Range("L3").Value = "Charisma=1/Thaumaturgy/0=3/Ray of Sickness/2=5/Crown of Madness/2"
$ws.Range("M3").Value = "Charisma=1/Thaumaturgy/0=3/Disguise Self/1=5/Detect Thoughts/2"
$ws.Range("N3").Value = "Charisma=1/Friends/0=3/Charm Person/2=5/Suggestion/2"
$ws.Range("O3").Value = "Charisma=1/Minor Illusion/0=3/Disguise Self/1=5/Invisibility/2"
$ws.Range("P3").Value = "Charisma=1/Ray of Frost/0=3/Armor of Agathys/2=5/Darkness/2"
$ws.Range("Q3").Value = "Charisma=1/Mage Hand/0=3/Tenser's Floating Disk/1=5/Arcane Lock/2"
$ws.Range("R3").Value = "Charisma=1/Mage Hand/0=3/Burning Hands/2=5/Flame Blade/2"
$ws.Range("S3").Value = "Charisma=1/Thaumaturgy/0=3/Searing Smite/2=5/Branding Smite/2"
$ws.Range("K4").Value = "Darkvision=Hellish Resistance"
$ws.Range("L4").Value = "Darkvision=Hellish Resistance"
$ws.Range("M4").Value = "Darkvision=Hellish Resistance"
$ws.Range("N4").Value = "Darkvision=Hellish Resistance"
$ws.Range("O4").Value = "Darkvision=Hellish Resistance"
$ws.Range("P4").Value = "Darkvision=Hellish Resistance"
$ws.Range("Q4").Value = "Darkvision=Hellish Resistance"
$ws.Range("R4").Value = "Darkvision=Hellish Resistance"
$ws.Range("S4").Value = "Darkvision=Hellish Resistance"
$ws.Range("K5").Value = 150.0
$ws.Range("L5").Value = 150.0
$ws.Range("M5").Value = 150.0
$ws.Range("N5").Value = 150.0
$ws.Range("O5").Value = 150.0
$ws.Range("P5").Value = 150.0
$ws.Range("Q5").Value = 150.0
$ws.Range("R5").Value = 150.0
$ws.Range("S5").Value = 150.0
$ws.Range("K6").Value = "Medium"
$ws.Range("L6").Value = "Medium"
$ws.Range("M6").Value = "Medium"
$ws.Range("N6").Value = "Medium"
$ws.Range("O6").Value = "Medium"
$ws.Range("P6").Value = "Medium"
$ws.Range("Q6").Value = "Medium"
$ws.Range("R6").Value = "Medium"
$ws.Range("S6").Value = "Medium"
$ws.Range("K7").Value = 30.0
$ws.Range("L7").Value = 30.0
$ws.Range("M7").Value = 30.0
$ws.Range("N7").Value = 30.0
$ws.Range("O7").Value = 30.0
$ws.Range("P7").Value = 30.0
$ws.Range("Q7").Value = 30.0
$ws.Range("R7").Value = 30.0
$ws.Range("S7").Value = 30.0
$ws.Range("K8").Value = "Common=Infernal"
$ws.Range("L8").Value = "Common=Infernal"
$ws.Range("M8").Value = "Common=Infernal"
$ws.Range("N8").Value = "Common=Infernal"
$ws.Range("O8").Value = "Common=Infernal"
$ws.Range("P8").Value = "Common=Infernal"
$ws.Range("Q8").Value = "Common=Infernal"
$ws.Range("R8").Value = "Common=Infernal"
$ws.Range("S8").Value = "Common=Infernal"
$ws.Range("K9").Value = "races/tiefling/Names.xlsx"
$ws.Range("L9").Value = "races/tiefling/Names.xlsx"
$ws.Range("M9").Value = "races/tiefling/Names.xlsx"
$ws.Range("N9").Value = "races/tiefling/Names.xlsx"
$ws.Range("O9").Value = "races/tiefling/Names.xlsx"
$ws.Range("P9").Value = "races/tiefling/Names.xlsx"
$ws.Range("Q9").Value = "races/tiefling/Names.xlsx"
$ws.Range("R9").Value = "races/tiefling/Names.xlsx"
$ws.Range("S9").Value = "races/tiefling/Names.xlsx"
